$d = $word.ActiveDocument

# 1. Fix the title typo "Лабораторная работная работа" -> "Лабораторная работа"
#    (drop the duplicated "работная ") and, as Word naturally does when the
#    cursor is left there after an edit, relocate the hidden "_GoBack" bookmark
#    to sit right after "Лабораторная ".
$d.Content.Find.Execute("Лабораторная работная работа", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Лабораторная работа", 2) | Out-Null

$titleRange = $d.Content
$titleRange.Find.Execute("Лабораторная ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$titleRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $titleRange) | Out-Null

# 2. Update the cached page-number field result in the footer from 7 to 8.
$footerRange = $d.Sections(1).Footers(1).Range
$pageNoChar = $footerRange.Characters(1)
$pageNoChar.Text = "8"
